$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.282"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05589"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.367"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8104"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9624"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1411"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07531"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03164"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03045"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09298"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.573"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001606"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04710"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005761"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006412"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005028"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001031"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.751"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.117"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3254"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1248"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'0.0003095"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03926"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007067"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1048"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003395"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008800"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005805"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.0005492"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.6790"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.1536"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01009"
$ws.Range("D51").Style = "Normal"
